$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Inclusão dos custos de transporte no orcamento
$ws.Range("A8").Value = "Transporte"
$ws.Range("B8").Value = "R$ 102"

# Column A was best-fit to the new (longer) content; Excel's ColumnWidth
# property applies a ~0.8333-character padding on top of the stored value,
# so back it out to land on the OOXML width of 44.
$ws.Columns.Item(1).ColumnWidth = 43.166666666666664

# Match the author's final selection / active cell.
$ws.Range("B8").Select()
